$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting the existing data
# (previously starting at column B) one column to the right.
$ws.Columns("B").Insert()

# The newly inserted column B was "best fit" sized before; give it the
# same visual width as column A/C but without the bestFit auto-size flag.
$ws.Columns("B").ColumnWidth = 14.6

# Populate the new column B with the "type" label for each data row.
$ws.Range("B9").Value = "EVP<2>"
$ws.Range("B10").Value = "Naility"
$ws.Range("B11").Value = "bool"
$ws.Range("B12").Value = "Naility"
$ws.Range("B13").Value = "Naility"
$ws.Range("B14").Value = "usize"

# Restore the user's on-screen selection to G3.
$ws.Range("G3").Select()
